# "adjusted weights to give better behaviour"
#
# The "For B2" weighting table (rows 27-32) originally reused the already
# weighted averages in D8/E8 ("avg velocity excluding b2") directly as the
# velocity-adjustment inputs in C29/D29. This introduces a new, explicit
# "Velocity adjustment" helper (D9 label, D10/E10 = D8-D3 / E8-E3, i.e. the
# excess of the flock's average velocity over b2's own velocity) and points
# C29/D29 at that instead, giving better (less runaway) behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label + helper formulas just below the existing "avg velocity
# excluding b2" row (row 8).
$ws.Range("D9").Value = "Velocity adjustment"
$ws.Range("D10").Formula = "=D8-D3"
$ws.Range("E10").Formula = "=E8-E3"

# Point the "For B2" velocity-adjustment inputs at the new helper cells
# instead of directly at D8/E8.
$ws.Range("C29").Formula = "=D10"
$ws.Range("D29").Formula = "=E10"

# Move the active selection to C32 (matches the author's saved cursor
# position after making the edit).
[void]$ws.Range("C32").Select()

[void]$wb.Application.Calculate()
